# Update the acquisition timestamp (取得日時) in column A of the
# "ランサーズ" sheet for the existing data rows (2-7) to reflect the
# new scrape run time: 2025-09-15 06:35:56 (appended at 2025-09-15 06:35 JST).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-15 06:35:56"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
